$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "Modelo", copying the formatting used by the other
# header cells (bold font + border + centered alignment) from E1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "Modelo"

# Update existing numeric predictions in row 2
$ws.Range("B2").Value = 0.3450020433581751
$ws.Range("C2").Value = 0.9933629278466505
$ws.Range("D2").Value = 0.4689456870504659

# Add new model-name cell F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"
